$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs_for_tollcalib")

# New tollclass segments added: tollclasses broken into 13 segments by direction
# (NGF_BPALT13Segments network project)
$names = @(
    "101_Marin_N - NGF_BPALT13Segments",
    "101_Marin_S - NGF_BPALT13Segments",
    "101_Peninsula_N - NGF_BPALT13Segments",
    "101_Peninsula_S - NGF_BPALT13Segments",
    "237_E - NGF_BPALT13Segments",
    "237_W - NGF_BPALT13Segments",
    "238_N - NGF_BPALT13Segments",
    "238_S - NGF_BPALT13Segments",
    "280_N - NGF_BPALT13Segments",
    "280_S - NGF_BPALT13Segments",
    "380_E - NGF_BPALT13Segments",
    "380_W - NGF_BPALT13Segments",
    "580_E - NGF_BPALT13Segments",
    "580_W - NGF_BPALT13Segments",
    "680_N - NGF_BPALT13Segments",
    "680_S - NGF_BPALT13Segments",
    "80_E - NGF_BPALT13Segments",
    "80_W - NGF_BPALT13Segments",
    "85_87_N - NGF_BPALT13Segments",
    "85_87_S - NGF_BPALT13Segments",
    "880_17_N - NGF_BPALT13Segments",
    "880_17_S - NGF_BPALT13Segments",
    "92_E - NGF_BPALT13Segments",
    "92_W - NGF_BPALT13Segments",
    "980_24_680_242_4_E - NGF_BPALT13Segments",
    "980_24_680_242_4_W - NGF_BPALT13Segments"
)

$startRow = 294
$startTollclass = 980001

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $startTollclass + $i
}

# Match the author's final selection/scroll position after adding the rows
$ws.Activate() | Out-Null
$ws.Range("A305").Select() | Out-Null

